# Insert a new data row at row 104 (pushes existing rows 104-222 down to 105-223,
# carrying their values/formatting with them - matches the diff's net effect of
# a weekly "prepend latest reading" update), then populate the new row with the
# newest reading's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 104..222 down to 105..223.
$ws.Rows.Item(104).Insert()

# Fill in the new row 104 with the latest observation. Columns A,B,C,E,F,G,H,I,N,Q,R
# are constant across every data row in this sheet (same market/category metadata),
# only D (fecha), J (volumen), K/L/M (precios), P (precio $/Kg) vary per-row.
$ws.Cells.Item(104, 1).Value = 7
$ws.Cells.Item(104, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(104, 3).Value = "Ñuble"
$ws.Cells.Item(104, 4).Value = 44546
$ws.Cells.Item(104, 5).Value = 16
$ws.Cells.Item(104, 6).Value = 100114013
$ws.Cells.Item(104, 7).Value = "Zanahoria"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 200
$ws.Cells.Item(104, 11).Value = 7000
$ws.Cells.Item(104, 12).Value = 7500
$ws.Cells.Item(104, 13).Value = 7250
$ws.Cells.Item(104, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(104, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(104, 16).Value = 362
$ws.Cells.Item(104, 17).Value = 20
$ws.Cells.Item(104, 18).Value = "Hortaliza"
